$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 182, shifting existing rows 182-191 down to 183-192.
$ws.Rows.Item(182).Insert()

# Populate the new row 182 with the new weekly data point.
$ws.Cells.Item(182, 1).Value = 11
$ws.Cells.Item(182, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(182, 3).Value = "Bíobío"
$ws.Cells.Item(182, 4).Value = 45041
$ws.Cells.Item(182, 5).Value = 8
$ws.Cells.Item(182, 6).Value = 100112021
$ws.Cells.Item(182, 7).Value = "Ají"
$ws.Cells.Item(182, 8).Value = "Americana (o)"
$ws.Cells.Item(182, 9).Value = "Primera"
$ws.Cells.Item(182, 10).Value = 25
$ws.Cells.Item(182, 11).Value = 17000
$ws.Cells.Item(182, 12).Value = 18000
$ws.Cells.Item(182, 13).Value = 17400
$ws.Cells.Item(182, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(182, 15).Value = "Región Metropolitana"
$ws.Cells.Item(182, 16).Value = 696
$ws.Cells.Item(182, 17).Value = 25
$ws.Cells.Item(182, 18).Value = "Hortaliza"
